$d = $word.ActiveDocument

# 1) Replace the long paragraph's trailing portion with the new wording
#    (generalized Levenshtein edit distance / threshold phrasing).
$old1 = "This was done by computing the edit distance (edits such as deletion, insertion or substitutions needed to convert one string to another string) of the clinical trial disease to each WHO database term and then if there were any matches within 20% any WHO tumor names then it was flagged as a potential tumor."
$new1 = "This was done by computing the generalized Levenshtein edit distance of the clinical trial disease to each WHO database term. If the clinical trial disease name was within a certain predetermined threshold (defined in algorithm 1) of the generalized Levenshtein edit distance  "

$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# 2) Find the paragraph that now ends with "...edit distance  " (the one we
#    just edited) and the "2.3 Adult and Pediatric Tumor Annotation of
#    Disease Data" heading. Between them sits a run of leftover scratch
#    paragraphs (blank lines, "As evident from table 2...", "WHO and NCIT
#    for standardizing", "Fuzzy match", "Manual annotation.", etc.) that
#    must be removed, keeping just a single paragraph for the new
#    continuation sentence.
$editedPara = $null
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($editedPara -eq $null -and $t -match "predetermined threshold") {
        $editedPara = $i
    }
    if ($t -match "2\.3 Adult and Pediatric Tumor Annotation") {
        $headingPara = $i
        break
    }
}

$startPara = $editedPara + 1
$pStart = $d.Paragraphs.Item($startPara)
$pHeading = $d.Paragraphs.Item($headingPara)

$deleteBegin = $pStart.Range.End
$deleteEnd = $pHeading.Range.Start

if ($deleteEnd -gt $deleteBegin) {
    $killRange = $d.Range($deleteBegin, $deleteEnd)
    $killRange.Delete()
}

# 3) Fill the now-solitary scratch paragraph with the new continuation text.
$new2 = "any WHO tumor name then it was flagged as a potential tumor. Once every disease in clinical trial was flagged as a potential tumor using steps 1 and 2 , then they were manually validated as tumors. Furthermore during the validation process we also annotated whether the disease was a pediatric tumor and added a field containing a citation that suggested that the tumor was a pediatric tumor. "

$pStart = $d.Paragraphs.Item($startPara)
$fillRange = $d.Range($pStart.Range.Start, $pStart.Range.End - 1)
$fillRange.Text = $new2
